{"js": "// Remove the trailing \"Ver no Jupiter ...\" / \"\u00a9 2020 ...\" footer block\n// (and the blank paragraph directly above it) that used to follow the\n// \"LOB1012: Estat\u00edstica (Requisito)\" paragraph at the end of the document.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the two footer paragraphs by their exact text.\nlet jupiterIdx = -1;\nlet copyrightIdx = -1;\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text.trim();\n  if (t === \"Ver no Jupiter Salvar em pdf Salvar em docx\") {\n    jupiterIdx = i;\n  } else if (t.startsWith(\"\u00a9 2020\")) {\n    copyrightIdx = i;\n  }\n}\n\nif (jupiterIdx !== -1 && copyrightIdx !== -1) {\n  // The blank paragraph immediately preceding the \"Ver no Jupiter\" line is\n  // part of the block being removed (it separated the requirements section\n  // from the footer).\n  let blankIdx = jupiterIdx - 1;\n  if (blankIdx >= 0 && items[blankIdx].text.trim() === \"\") {\n    items[blankIdx].delete();\n  }\n  items[jupiterIdx].delete();\n  items[copyrightIdx].delete();\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter ...\" / \"\u00a9 2020 ...\" footer block\n# (and the blank paragraph directly above it) that used to follow the\n# \"LOB1012: Estat\u00edstica (Requisito)\" paragraph at the end of the document.\n\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n$jupiterIdx = -1\n$copyrightIdx = -1\n\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.Trim()\n    if ($t -eq \"Ver no Jupiter Salvar em pdf Salvar em docx\") {\n        $jupiterIdx = $i\n    } elseif ($t.StartsWith(\"\u00a9 2020\")) {\n        $copyrightIdx = $i\n    }\n}\n\nif ($jupiterIdx -ne -1 -and $copyrightIdx -ne -1) {\n    # The blank paragraph immediately preceding the \"Ver no Jupiter\" line is\n    # part of the block being removed (it separated the requirements\n    # section from the footer).\n    $blankIdx = $jupiterIdx - 1\n    $hasBlank = $false\n    if ($blankIdx -ge 1) {\n        $blankText = $d.Paragraphs.Item($blankIdx).Range.Text.Trim()\n        if ($blankText -eq \"\") {\n            $hasBlank = $true\n        }\n    }\n\n    # Delete from the highest index down to the lowest so earlier indices\n    # stay valid while later ones are removed.\n    $d.Paragraphs.Item($copyrightIdx).Range.Delete()\n    $d.Paragraphs.Item($jupiterIdx).Range.Delete()\n    if ($hasBlank) {\n        $d.Paragraphs.Item($blankIdx).Range.Delete()\n    }\n}\n"}
